$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3395
$ws.Range("J17").Value = 3774
$ws.Range("L17").Value = 11322
$ws.Range("N17").Value = -11658
$ws.Range("H28").Value = 964.0789
$ws.Range("I28").Value = 824.75
$ws.Range("K28").Value = 824.75
$ws.Range("M28").Value = -339.75
$ws.Range("H86").Value = 6583672
$ws.Range("I86").Value = 7400
$ws.Range("K86").Value = 7400
$ws.Range("M86").Value = -6277
$ws.Range("H89").Value = 6583672
$ws.Range("I89").Value = 7400
$ws.Range("K89").Value = 37000
$ws.Range("M89").Value = -31384
$ws.Range("H132").Value = 23644.738
$ws.Range("I132").Value = 1780.4736
$ws.Range("K132").Value = 5341.4208
$ws.Range("M132").Value = -2811.4208
$ws.Range("H138").Value = 1444.5927
$ws.Range("I138").Value = 971.1905
$ws.Range("J138").Value = 3101.5
$ws.Range("K138").Value = 2913.5715
$ws.Range("L138").Value = 9304.5
$ws.Range("M138").Value = 2226.4285
$ws.Range("N138").Value = -19584.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 212
$ws.Range("I5").Value = 156.66667
$ws.Range("J5").Value = 267.33334
$ws.Range("K5").Value = 156.66667
$ws.Range("L5").Value = 267.33334
$ws.Range("M5").Value = -44.66667000000001
$ws.Range("N5").Value = -491.33334
$ws.Range("H97").Value = 3832237.5
$ws.Range("I97").Value = 871.6818
$ws.Range("K97").Value = 871.6818
$ws.Range("M97").Value = -375.6818

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 212
$ws.Range("I4").Value = 156.66667
$ws.Range("J4").Value = 267.33334
$ws.Range("K4").Value = 156.66667
$ws.Range("L4").Value = 267.33334
$ws.Range("M4").Value = -41.66667000000001
$ws.Range("N4").Value = -497.33334
$ws.Range("H22").Value = 326.66666
$ws.Range("I22").Value = 352.5
$ws.Range("J22").Value = 275
$ws.Range("K22").Value = 352.5
$ws.Range("L22").Value = 275
$ws.Range("M22").Value = -179.5
$ws.Range("N22").Value = -621
$ws.Range("H94").Value = 5955717.5
$ws.Range("I94").Value = 2918.6667
$ws.Range("J94").Value = 12824332
$ws.Range("K94").Value = 2918.6667
$ws.Range("L94").Value = 12824332
$ws.Range("M94").Value = -2467.6667
$ws.Range("N94").Value = -12825234
$ws.Range("H105").Value = 4749.048
$ws.Range("I105").Value = 4291.4287
$ws.Range("J105").Value = 5664.2856
$ws.Range("K105").Value = 4291.4287
$ws.Range("L105").Value = 5664.2856
$ws.Range("M105").Value = -2544.4287
$ws.Range("N105").Value = -9158.285599999999
$ws.Range("H134").Value = 2693.6738
$ws.Range("I134").Value = 1703.8108
$ws.Range("J134").Value = 6763.1113
$ws.Range("K134").Value = 5111.4324
$ws.Range("L134").Value = 20289.3339
$ws.Range("M134").Value = -2576.4324
$ws.Range("N134").Value = -25359.3339
$ws.Range("H138").Value = 78534.8
$ws.Range("J138").Value = 78697.64
$ws.Range("L138").Value = 78697.64
$ws.Range("N138").Value = -88977.64

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10206055
$ws.Range("I99").Value = 1747243.2
$ws.Range("K99").Value = 1747243.2
$ws.Range("M99").Value = -1745745.2
$ws.Range("H126").Value = 10206055
$ws.Range("I126").Value = 1747243.2
$ws.Range("K126").Value = 5241729.6
$ws.Range("M126").Value = -5239259.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 211.18182
$ws.Range("J107").Value = 217.3
$ws.Range("L107").Value = 651.9000000000001
$ws.Range("N107").Value = -4491.9
$ws.Range("H113").Value = 1134
$ws.Range("J113").Value = 1011.7857
$ws.Range("L113").Value = 3035.3571
$ws.Range("N113").Value = -7375.3571
$ws.Range("H121").Value = 2636.889
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2636.889
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 7910.667
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -10530.667
$ws.Range("H131").Value = 6483.125
$ws.Range("I131").Value = 4159.8
$ws.Range("J131").Value = 7539.1816
$ws.Range("K131").Value = 12479.4
$ws.Range("L131").Value = 22617.5448
$ws.Range("M131").Value = -7439.400000000001
$ws.Range("N131").Value = -32697.5448
$ws.Range("H132").Value = 2239
$ws.Range("I132").Value = 2173.75
$ws.Range("K132").Value = 19563.75
$ws.Range("M132").Value = -17033.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 501.5
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H20").Value = 14004699
$ws.Range("I20").Value = 17500000
$ws.Range("J20").Value = 23495
$ws.Range("K20").Value = 17500000
$ws.Range("L20").Value = 23495
$ws.Range("M20").Value = -17499755
$ws.Range("N20").Value = -23985
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H102").Value = 2073.926
$ws.Range("I102").Value = 1284.1578
$ws.Range("J102").Value = 3949.625
$ws.Range("K102").Value = 1284.1578
$ws.Range("L102").Value = 3949.625
$ws.Range("M102").Value = 337.8422
$ws.Range("N102").Value = -7193.625
$ws.Range("H122").Value = 3236.5386
$ws.Range("I122").Value = 1244.1765
$ws.Range("J122").Value = 6999.8887
$ws.Range("K122").Value = 3732.5295
$ws.Range("L122").Value = 20999.6661
$ws.Range("M122").Value = -1282.5295
$ws.Range("N122").Value = -25899.6661

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1139.5
$ws.Range("I16").Value = 385.75
$ws.Range("J16").Value = 2647
$ws.Range("K16").Value = 385.75
$ws.Range("L16").Value = 2647
$ws.Range("M16").Value = -215.75
$ws.Range("N16").Value = -2987
$ws.Range("H46").Value = 6396.52
$ws.Range("I46").Value = 673.4286
$ws.Range("J46").Value = 8622.166999999999
$ws.Range("K46").Value = 673.4286
$ws.Range("L46").Value = 8622.166999999999
$ws.Range("M46").Value = -485.4286
$ws.Range("N46").Value = -8998.166999999999
$ws.Range("H122").Value = 6030.3335
$ws.Range("J122").Value = 7105.8887
$ws.Range("L122").Value = 21317.6661
$ws.Range("N122").Value = -26217.6661

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2667.318
$ws.Range("J132").Value = 3542
$ws.Range("L132").Value = 10626
$ws.Range("N132").Value = -15686
